{"js": "// Insert the \"No report for week ending 3/17/2019 ...\" paragraph content\n// into the existing empty paragraph that currently only holds the\n// `_GoBack` bookmark. Locate that exact spot via the bookmark itself\n// (robust against paragraph-index/ordering assumptions) rather than a\n// hard-coded paragraph index.\nconst body = context.document.body;\nconst target = body.getBookmarkRangeOrNullObject(\"_GoBack\");\ntarget.load(\"isNullObject\");\nawait context.sync();\n\nif (target.isNullObject) {\n  throw new Error(\"_GoBack bookmark not found in document body\");\n}\n\n// The OOXML package fragment containing the runs (plus the proofErr\n// spell/grammar-check markers present in the original authored content)\n// to insert at the very start of that paragraph, ahead of the bookmark.\nconst ooxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:r><w:t>No report for week e</w:t></w:r>' +\n              '<w:r><w:t>nding 3/17</w:t></w:r>' +\n              '<w:r><w:t>/2019 \\u2013 Spring break.</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> Did an individual code review to suggest bug fixes and implementation of the </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:proofErr w:type=\"gramStart\"/>' +\n              '<w:r><w:t>addTeam</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t>(</w:t></w:r>' +\n              '<w:proofErr w:type=\"gramEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\">) and </w:t></w:r>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r><w:t>removeTeam</w:t></w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\">() functions. Some of the bug fixes are simple such as the use of a while loop instead of </w:t></w:r>' +\n              '<w:proofErr w:type=\"gramStart\"/>' +\n              '<w:r><w:t>a for</w:t></w:r>' +\n              '<w:proofErr w:type=\"gramEnd\"/>' +\n              '<w:r><w:t xml:space=\"preserve\"> loop. </w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">These items will be </w:t></w:r>' +\n              '<w:r><w:t>discussed</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> during our team meeting on 3/19/2019.</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, \"Start\");\nawait context.sync();\n\n", "ps1": "# Insert the \"No report for week ending 3/17/2019 ...\" paragraph content\n# into the existing (empty) paragraph that currently only holds the\n# `_GoBack` bookmark. The new runs land immediately before the bookmark,\n# inside the same paragraph, exactly like the authored edit.\n$d = $word.ActiveDocument\n\n# Locate the _GoBack bookmark and build a zero-length (collapsed) Range at\n# its start so the inserted XML lands right before <w:bookmarkStart>.\n$bmStart = $d.Bookmarks.Item(\"_GoBack\").Range.Start\n$r = $d.Range($bmStart, $bmStart)\n\n# Flat-OPC package fragment with the new runs, including the proofErr\n# spell-check/grammar-check markers present in the authored content\n# (cosmetic markers around \"addTeam(\", \"removeTeam()\" and \"a for\").\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>No report for week e</w:t></w:r><w:r><w:t>nding 3/17</w:t></w:r><w:r><w:t>/2019 \u2013 Spring break.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Did an individual code review to suggest bug fixes and implementation of the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>addTeam</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\">) and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>removeTeam</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">() functions. Some of the bug fixes are simple such as the use of a while loop instead of </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>a for</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> loop. </w:t></w:r><w:r><w:t xml:space=\"preserve\">These items will be </w:t></w:r><w:r><w:t>discussed</w:t></w:r><w:r><w:t xml:space=\"preserve\"> during our team meeting on 3/19/2019.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($xml)\n"}
